$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is empty for every row and was removed
# from the export, shifting the subsequent columns (reviews_average,
# latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) one position to the left.
$ws.Columns.Item(5).Delete()
